$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A11").Value = 10
$ws.Range("B11").Value = "hi manager how are you"
